$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 50.80855966666667
$ws.Range("H2").Value = 152.425679
$ws.Range("I2").Value = 0.1328146143749533
$ws.Range("J2").Value = 0.1328146143749533
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 6.840255666666667
$ws.Range("N2").Value = 20.520767
$ws.Range("O2").Value = 0.2326221245729483
$ws.Range("P2").Value = 0.2326221245729483
$ws.Range("Q2").Value = 347.5435381750881
$ws.Range("R2").Value = 3127.891843575793
$ws.Range("S2").Value = 0.03089561777023849
$ws.Range("T2").Value = 0.03089561777023849

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 50.80855966666667
$ws.Range("H3").Value = 152.425679
$ws.Range("I3").Value = 0.1328146143749533
$ws.Range("J3").Value = 0.1328146143749533
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 9.415322999999999
$ws.Range("N3").Value = 28.245969
$ws.Range("O3").Value = 0.3201945287620894
$ws.Range("P3").Value = 0.3201945287620895
$ws.Range("Q3").Value = 478.3790004264389
$ws.Range("R3").Value = 4305.411003837951
$ws.Range("S3").Value = 0.04252651286250681
$ws.Range("T3").Value = 0.04252651286250682

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 50.80855966666667
$ws.Range("H4").Value = 152.425679
$ws.Range("I4").Value = 0.1328146143749533
$ws.Range("J4").Value = 0.1328146143749533
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 7.720817
$ws.Range("N4").Value = 23.162451
$ws.Range("O4").Value = 0.2625680883144773
$ws.Range("P4").Value = 0.2625680883144773
$ws.Range("Q4").Value = 392.2835912199143
$ws.Range("R4").Value = 3530.552320979229
$ws.Range("S4").Value = 0.034872879396656
$ws.Range("T4").Value = 0.034872879396656

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 50.80855966666667
$ws.Range("H5").Value = 152.425679
$ws.Range("I5").Value = 0.1328146143749533
$ws.Range("J5").Value = 0.1328146143749533
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 5.428613333333334
$ws.Range("N5").Value = 16.28584
$ws.Range("O5").Value = 0.1846152583504849
$ws.Range("P5").Value = 0.1846152583504849
$ws.Range("Q5").Value = 275.8200244539289
$ws.Range("R5").Value = 2482.38022008536
$ws.Range("S5").Value = 0.02451960434555204
$ws.Range("T5").Value = 0.02451960434555204

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 6.697301333333333
$ws.Range("H6").Value = 20.091904
$ws.Range("I6").Value = 0.01750688269408059
$ws.Range("J6").Value = 0.01750688269408059
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 6.840255666666667
$ws.Range("N6").Value = 20.520767
$ws.Range("O6").Value = 0.2326221245729483
$ws.Range("P6").Value = 0.2326221245729483
$ws.Range("Q6").Value = 45.81125339670756
$ws.Range("R6").Value = 412.301280570368
$ws.Range("S6").Value = 0.004072488246946408
$ws.Range("T6").Value = 0.004072488246946408

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 6.697301333333333
$ws.Range("H7").Value = 20.091904
$ws.Range("I7").Value = 0.01750688269408059
$ws.Range("J7").Value = 0.01750688269408059
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 9.415322999999999
$ws.Range("N7").Value = 28.245969
$ws.Range("O7").Value = 0.3201945287620894
$ws.Range("P7").Value = 0.3201945287620895
$ws.Range("Q7").Value = 63.05725528166399
$ws.Range("R7").Value = 567.515297534976
$ws.Range("S7").Value = 0.005605608054324312
$ws.Range("T7").Value = 0.005605608054324313

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 6.697301333333333
$ws.Range("H8").Value = 20.091904
$ws.Range("I8").Value = 0.01750688269408059
$ws.Range("J8").Value = 0.01750688269408059
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 7.720817
$ws.Range("N8").Value = 23.162451
$ws.Range("O8").Value = 0.2625680883144773
$ws.Range("P8").Value = 0.2625680883144773
$ws.Range("Q8").Value = 51.70863798852267
$ws.Range("R8").Value = 465.377741896704
$ws.Range("S8").Value = 0.004596748721330546
$ws.Range("T8").Value = 0.004596748721330546

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 6.697301333333333
$ws.Range("H9").Value = 20.091904
$ws.Range("I9").Value = 0.01750688269408059
$ws.Range("J9").Value = 0.01750688269408059
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 5.428613333333334
$ws.Range("N9").Value = 16.28584
$ws.Range("O9").Value = 0.1846152583504849
$ws.Range("P9").Value = 0.1846152583504849
$ws.Range("Q9").Value = 36.35705931548444
$ws.Range("R9").Value = 327.21353383936
$ws.Range("S9").Value = 0.003232037671479321
$ws.Range("T9").Value = 0.003232037671479321

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 311.72598
$ws.Range("H10").Value = 935.1779399999999
$ws.Range("I10").Value = 0.8148580887939706
$ws.Range("J10").Value = 0.8148580887939705
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 6.840255666666667
$ws.Range("N10").Value = 20.520767
$ws.Range("O10").Value = 0.2326221245729483
$ws.Range("P10").Value = 0.2326221245729483
$ws.Range("Q10").Value = 2132.28540114222
$ws.Range("R10").Value = 19190.56861027998
$ws.Range("S10").Value = 0.1895540198407056
$ws.Range("T10").Value = 0.1895540198407056

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 311.72598
$ws.Range("H11").Value = 935.1779399999999
$ws.Range("I11").Value = 0.8148580887939706
$ws.Range("J11").Value = 0.8148580887939705
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 9.415322999999999
$ws.Range("N11").Value = 28.245969
$ws.Range("O11").Value = 0.3201945287620894
$ws.Range("P11").Value = 0.3201945287620895
$ws.Range("Q11").Value = 2935.00078919154
$ws.Range("R11").Value = 26415.00710272386
$ws.Range("S11").Value = 0.2609131017493622
$ws.Range("T11").Value = 0.2609131017493622

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 311.72598
$ws.Range("H12").Value = 935.1779399999999
$ws.Range("I12").Value = 0.8148580887939706
$ws.Range("J12").Value = 0.8148580887939705
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 7.720817
$ws.Range("N12").Value = 23.162451
$ws.Range("O12").Value = 0.2625680883144773
$ws.Range("P12").Value = 0.2625680883144773
$ws.Range("Q12").Value = 2406.77924572566
$ws.Range("R12").Value = 21661.01321153094
$ws.Range("S12").Value = 0.2139557306222215
$ws.Range("T12").Value = 0.2139557306222214

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 311.72598
$ws.Range("H13").Value = 935.1779399999999
$ws.Range("I13").Value = 0.8148580887939706
$ws.Range("J13").Value = 0.8148580887939705
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 5.428613333333334
$ws.Range("N13").Value = 16.28584
$ws.Range("O13").Value = 0.1846152583504849
$ws.Range("P13").Value = 0.1846152583504849
$ws.Range("Q13").Value = 1692.2398113744
$ws.Range("R13").Value = 15230.1583023696
$ws.Range("S13").Value = 0.1504352365816812
$ws.Range("T13").Value = 0.1504352365816812

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 13.32063566666667
$ws.Range("H14").Value = 39.961907
$ws.Range("I14").Value = 0.03482041413699558
$ws.Range("J14").Value = 0.03482041413699557
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 6.840255666666667
$ws.Range("N14").Value = 20.520767
$ws.Range("O14").Value = 0.2326221245729483
$ws.Range("P14").Value = 0.2326221245729483
$ws.Range("Q14").Value = 91.11655360251879
$ws.Range("R14").Value = 820.048982422669
$ws.Range("S14").Value = 0.008099998715057836
$ws.Range("T14").Value = 0.008099998715057835

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 13.32063566666667
$ws.Range("H15").Value = 39.961907
$ws.Range("I15").Value = 0.03482041413699558
$ws.Range("J15").Value = 0.03482041413699557
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 9.415322999999999
$ws.Range("N15").Value = 28.245969
$ws.Range("O15").Value = 0.3201945287620894
$ws.Range("P15").Value = 0.3201945287620895
$ws.Range("Q15").Value = 125.418087366987
$ws.Range("R15").Value = 1128.762786302883
$ws.Range("S15").Value = 0.01114930609589609
$ws.Range("T15").Value = 0.01114930609589609

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 13.32063566666667
$ws.Range("H16").Value = 39.961907
$ws.Range("I16").Value = 0.03482041413699558
$ws.Range("J16").Value = 0.03482041413699557
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 7.720817
$ws.Range("N16").Value = 23.162451
$ws.Range("O16").Value = 0.2625680883144773
$ws.Range("P16").Value = 0.2625680883144773
$ws.Range("Q16").Value = 102.8461903060063
$ws.Range("R16").Value = 925.6157127540571
$ws.Range("S16").Value = 0.009142729574269329
$ws.Range("T16").Value = 0.009142729574269327

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 13.32063566666667
$ws.Range("H17").Value = 39.961907
$ws.Range("I17").Value = 0.03482041413699558
$ws.Range("J17").Value = 0.03482041413699557
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 5.428613333333334
$ws.Range("N17").Value = 16.28584
$ws.Range("O17").Value = 0.1846152583504849
$ws.Range("P17").Value = 0.1846152583504849
$ws.Range("Q17").Value = 72.31258038854223
$ws.Range("R17").Value = 650.81322349688
$ws.Range("S17").Value = 0.006428379751772315
$ws.Range("T17").Value = 0.006428379751772314
